$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Updated source values (re-run of the comparison with new inputs)
#    NB: the PowerShell parser here does not accept scientific-notation
#    numeric literals (e.g. "9.8E-05"), so every value below is written out
#    in plain decimal form (with enough digits to round-trip to the exact
#    same IEEE-754 double as the target value).
# ---------------------------------------------------------------------------
$ws.Range("C5").Value = 0.0000979999999999999968359

$ws.Range("B8").Value = 0.366005
$ws.Range("C8").Value = 0.370227
$ws.Range("D8").Value = 0.048825

$ws.Range("B9").Value = 0.043282
$ws.Range("C9").Value = 0.04325
$ws.Range("D9").Value = 0.005765

$ws.Range("B10").Value = 0.003142
$ws.Range("C10").Value = 0.003134
$ws.Range("D10").Value = 0.000425

$ws.Range("B13").Value = 0.161871985
$ws.Range("C13").Value = 0.000991949
$ws.Range("D13").Value = 0.0000525060000000000017808

$ws.Range("B14").Value = 0.002263353
$ws.Range("C14").Value = 0.0000126830000000000003716
$ws.Range("D14").Value = 0.0000007369999999999999952

$ws.Range("B15").Value = 0.0000132000000000000006732
$ws.Range("C15").Value = 0.0000000800000000000000017
$ws.Range("D15").Value = 0.0000000200000000000000004
$ws.Range("F15").Value = 0.0000000728000000000000034
$ws.Range("F15").NumberFormat = "0.00E+00"
$ws.Range("F15").Font.Color = $ws.Range("F9").Font.Color
$ws.Range("G15").Formula = "=0.00000000402"

$ws.Range("B17").Value = 0.4135972
$ws.Range("C17").Value = 0.09025124
$ws.Range("D17").Value = 0.00789718

$ws.Range("B18").Value = 0.183851182
$ws.Range("C18").Value = 0.047438268
$ws.Range("D18").Value = 0.00332179

$ws.Range("B19").Value = 0.066183403
$ws.Range("C19").Value = 0.016484322
$ws.Range("D19").Value = 0.001183312

$ws.Range("B25").Value = 0.161847
$ws.Range("C25").Value = 0.001009
$ws.Range("D25").Value = 0.0000689999999999999966659

$ws.Range("B26").Value = 0.002267
$ws.Range("C26").Value = 0.002267
$ws.Range("D26").Value = 0.000361

# ---------------------------------------------------------------------------
# 2. Percent-format the "Mobley s/mean" columns (I:K), including the blank
#    separator rows, to match the new shared look of the comparison columns.
# ---------------------------------------------------------------------------
$ws.Range("I4:K26").NumberFormat = "0.00%"

# ---------------------------------------------------------------------------
# 3. New "AMC" / "Relative difference (%)" block in columns L:O
# ---------------------------------------------------------------------------
$ws.Range("L1").Value = "AMC"

$ws.Range("M1").Value = "Relative difference (%)"
$ws.Range("M1").NumberFormat = "0.00%"

$ws.Range("M2").Value = "Ed"
$ws.Range("N2").Value = "Eou"
$ws.Range("O2").Value = "Lu"
$ws.Range("M2:O2").NumberFormat = "0.00%"

$ws.Range("M4").Formula = "=(E4-B4)/E4"
$ws.Range("N4").Formula = "=(F4-C4)/F4"
$ws.Range("O4").Formula = "=(G4-D4)/G4"
$ws.Range("M5").Formula = "=(E5-B5)/E5"
$ws.Range("N5").Formula = "=(F5-C5)/F5"
$ws.Range("O5").Formula = "=(G5-D5)/G5"
$ws.Range("M6").Formula = "=(E6-B6)/E6"
$ws.Range("N6").Formula = "=(F6-C6)/F6"
$ws.Range("O6").Formula = "=(G6-D6)/G6"
$ws.Range("M4:O6").NumberFormat = "0.00%"

$ws.Range("M8").Formula = "=(E8-B8)/E8"
$ws.Range("N8").Formula = "=(F8-C8)/F8"
$ws.Range("O8").Formula = "=(G8-D8)/G8"
$ws.Range("M9").Formula = "=(E9-B9)/E9"
$ws.Range("N9").Formula = "=(F9-C9)/F9"
$ws.Range("O9").Formula = "=(G9-D9)/G9"
$ws.Range("M10").Formula = "=(E10-B10)/E10"
$ws.Range("N10").Formula = "=(F10-C10)/F10"
$ws.Range("O10").Formula = "=(G10-D10)/G10"
$ws.Range("M8:O10").NumberFormat = "0.00%"

$ws.Range("M13").Formula = "=(E13-B13)/E13"
$ws.Range("N13").Formula = "=(F13-C13)/F13"
$ws.Range("O13").Formula = "=(G13-D13)/G13"
$ws.Range("M14").Formula = "=(E14-B14)/E14"
$ws.Range("N14").Formula = "=(F14-C14)/F14"
$ws.Range("O14").Formula = "=(G14-D14)/G14"
$ws.Range("M15").Formula = "=(E15-B15)/E15"
$ws.Range("N15").Formula = "=(F15-C15)/F15"
$ws.Range("O15").Formula = "=(G15-D15)/G15"
$ws.Range("M13:O15").NumberFormat = "0.00%"

$ws.Range("M17").Formula = "=(E17-B17)/E17"
$ws.Range("N17").Formula = "=(F17-C17)/F17"
$ws.Range("O17").Formula = "=(G17-D17)/G17"
$ws.Range("M18").Formula = "=(E18-B18)/E18"
$ws.Range("N18").Formula = "=(F18-C18)/F18"
$ws.Range("O18").Formula = "=(G18-D18)/G18"
$ws.Range("M19").Formula = "=(E19-B19)/E19"
$ws.Range("N19").Formula = "=(F19-C19)/F19"
$ws.Range("O19").Formula = "=(G19-D19)/G19"
$ws.Range("M17:O19").NumberFormat = "0.00%"

$ws.Range("M25").Formula = "=(E25-B25)/E25"
$ws.Range("N25").Formula = "=(F25-C25)/F25"
$ws.Range("O25").Formula = "=(G25-D25)/G25"
$ws.Range("M26").Formula = "=(E26-B26)/E26"
$ws.Range("N26").Formula = "=(F26-C26)/F26"
$ws.Range("O26").Formula = "=(G26-D26)/G26"
$ws.Range("M25:O26").NumberFormat = "0.00%"

# ---------------------------------------------------------------------------
# 4. Column widths for the new block
# ---------------------------------------------------------------------------
$ws.Range("M:N").ColumnWidth = 9.125
$ws.Range("O:O").ColumnWidth = 9.625

# ---------------------------------------------------------------------------
# 5. Freeze the header rows and scroll/select like the saved view
# ---------------------------------------------------------------------------
$ws.Range("A3").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D17").Select()

# ---------------------------------------------------------------------------
# 6. Page setup
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1
